# Update crypto price/volume data per GitHub Actions scheduled refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.239.54'
$ws.Range("E2").Value = '  +0.06%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.596.11'
$ws.Range("E3").Value = '  +0.48%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.32'
$ws.Range("E5").Value = '  -0.22%  '
$ws.Range("E6").Value = '  +0.22%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  +0.23%  '
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.97'
$ws.Range("E10").Value = '  -1.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0855'
$ws.Range("E11").Value = '  +1.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.822.32'
$ws.Range("E12").Value = '  +0.57%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.607.84'
$ws.Range("E13").Value = '  +1.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.99'
$ws.Range("E14").Value = '  -0.48%  '
$ws.Range("E15").Value = '  -2.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.59'
$ws.Range("E16").Value = '  -0.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.258.79'
$ws.Range("E17").Value = '  +0.09%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '230.57'
$ws.Range("E18").Value = '  +7.75%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.66'
$ws.Range("E19").Value = '  +2.88%  '
$ws.Range("E20").Value = '  -0.41%  '
$ws.Range("E22").Value = '  -0.39%  '
$ws.Range("B23").Value = 'Toncoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.15'
$ws.Range("E23").Value = '  +1.72%  '
$ws.Range("B24").Value = 'Avalanche'
$ws.Range("C24").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.93'
$ws.Range("E24").Value = '  -0.50%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.16'
$ws.Range("E25").Value = '  +1.08%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.00'
$ws.Range("E27").Value = '  +0.32%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.32'
$ws.Range("E29").Value = '  +1.40%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0494'
$ws.Range("E30").Value = '  -0.18%  '
$ws.Range("E31").Value = '  -0.45%  '
$ws.Range("B32").Value = 'Maker'
$ws.Range("C32").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.475.90'
$ws.Range("E32").Value = '  +4.08%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.19'
$ws.Range("E33").Value = '  +0.41%  '
$ws.Range("E34").Value = '  -0.10%  '
$ws.Range("E35").Value = '  -0.73%  '
$ws.Range("E36").Value = '  +0.46%  '
$ws.Range("E37").Value = '  -3.86%  '
$ws.Range("E38").Value = '  -1.00%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.820'
$ws.Range("E39").Value = '  -0.16%  '
$ws.Range("E41").Value = '  +0.06%  '
$ws.Range("E42").Value = '  +1.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.934'
$ws.Range("E43").Value = '  -0.25%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.735.96'
$ws.Range("E44").Value = '  +0.73%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.755'
$ws.Range("E45").Value = '  -1.24%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.50'
$ws.Range("E46").Value = '  -1.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '87.90'
$ws.Range("E47").Value = '  +2.36%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.47'
$ws.Range("E48").Value = '  -0.91%  '
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("E50").Value = '  -1.92%  '
$ws.Range("E51").Value = '  -0.04%  '
